# Update NATMI Ltf-Lrp11 LR-pair sheet with newly recomputed TPM-derived stats.
# The "Sending cluster" for every data row changes from "Resolving-Mac" to "ECs",
# and the per-row expression / specificity metrics (columns E-T) are refreshed
# with the new TPM-based values. "Target cluster" (column D) values stay the
# same per row (ECs, FAPs, MuSCs, Resolving-Mac) and columns B/C (Ltf, Lrp11)
# are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("A2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03556333333333333
$ws.Range("H2").Value = 0.10669
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8541956666666667
$ws.Range("N2").Value = 2.562587
$ws.Range("O2").Value = 0.2037115339104792
$ws.Range("P2").Value = 0.2037115339104792
$ws.Range("Q2").Value = 0.03037804522555555
$ws.Range("R2").Value = 0.27340240703
$ws.Range("S2").Value = 0.2037115339104792
$ws.Range("T2").Value = 0.2037115339104792

# Row 3 (Target cluster: FAPs)
$ws.Range("A3").Value = "ECs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03556333333333333
$ws.Range("H3").Value = 0.10669
$ws.Range("O3").Value = 0.3552312021577347
$ws.Range("P3").Value = 0.3552312021577347
$ws.Range("Q3").Value = 0.05297309051444444
$ws.Range("R3").Value = 0.4767578146299999
$ws.Range("S3").Value = 0.3552312021577347
$ws.Range("T3").Value = 0.3552312021577347

# Row 4 (Target cluster: MuSCs)
$ws.Range("A4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03556333333333333
$ws.Range("H4").Value = 0.10669
$ws.Range("M4").Value = 1.431300333333333
$ws.Range("N4").Value = 4.293901
$ws.Range("O4").Value = 0.3413414487663211
$ws.Range("P4").Value = 0.3413414487663211
$ws.Range("Q4").Value = 0.05090181085444444
$ws.Range("R4").Value = 0.45811629769
$ws.Range("S4").Value = 0.3413414487663211
$ws.Range("T4").Value = 0.3413414487663211

# Row 5 (Target cluster: Resolving-Mac)
$ws.Range("A5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03556333333333333
$ws.Range("H5").Value = 0.10669
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.4181246666666666
$ws.Range("N5").Value = 1.254374
$ws.Range("O5").Value = 0.09971581516546497
$ws.Range("P5").Value = 0.09971581516546497
$ws.Range("Q5").Value = 0.01486990689555555
$ws.Range("R5").Value = 0.13382916206
$ws.Range("S5").Value = 0.09971581516546497
$ws.Range("T5").Value = 0.09971581516546497
